# Horarios actualizados Linea 141 - scrape refresh 08:34:05 -> 08:48:01
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with the
# newest scrape snapshot: refreshed "Hora_Scrap"/"Minutos" countdowns for
# rows whose scrape timestamp moved from 08:34:05 to 08:48:01, a handful of
# same-arrival-time rows that got reordered (their Linea/Hora_Scrap/Minutos
# values swapped places), and new arrival rows appended at the bottom of the
# LP1912 and 6203-6173 sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# ----- Sheet1: LP1912 -----
$ws1.Range("A2").Value = "Última actualización: 08:48:01"
$ws1.Range("A3").Value = "Total filas: 128"
$ws1.Range("C39").Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Range("C40").Value = "23_HERNANDEZ"
$ws1.Range("C47").Value = "23_HERNANDEZ"
$ws1.Range("C48").Value = "14_ABASTO"
$ws1.Range("A52").Value = "06:59:37"
$ws1.Range("C52").Value = "15_ABASTO"
$ws1.Range("D52").Value = 6
$ws1.Range("A53").Value = "05:52:07"
$ws1.Range("C53").Value = "23_HERNANDEZ"
$ws1.Range("D53").Value = 73
$ws1.Range("C63").Value = "11_ETCHEVERRY"
$ws1.Range("C64").Value = "16_SANTA ANA"
$ws1.Range("A65").Value = "05:52:07"
$ws1.Range("C65").Value = "11_ETCHEVERRY"
$ws1.Range("D65").Value = 100
$ws1.Range("C66").Value = "16_SANTA ANA"
$ws1.Range("A67").Value = "07:28:14"
$ws1.Range("C67").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D67").Value = 4
$ws1.Range("A72").Value = "07:28:14"
$ws1.Range("C72").Value = "16_SANTA ANA"
$ws1.Range("D72").Value = 19
$ws1.Range("A73").Value = "06:59:37"
$ws1.Range("C73").Value = "14_ABASTO"
$ws1.Range("D73").Value = 48
$ws1.Range("A88").Value = "07:28:14"
$ws1.Range("C88").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D88").Value = 55
$ws1.Range("A89").Value = "08:13:38"
$ws1.Range("C89").Value = "215B_EL PATO"
$ws1.Range("D89").Value = 10
$ws1.Range("A97").Value = "08:48:01"
$ws1.Range("D97").Value = 5
$ws1.Range("A101").Value = "08:48:01"
$ws1.Range("D101").Value = 14
$ws1.Range("A103").Value = "08:48:01"
$ws1.Range("D103").Value = 16
$ws1.Range("A104").Value = "08:48:01"
$ws1.Range("D104").Value = 17
$ws1.Range("A106").Value = "08:48:01"
$ws1.Range("D106").Value = 23
$ws1.Range("A107").Value = "08:48:01"
$ws1.Range("D107").Value = 25
$ws1.Range("A109").Value = "08:48:01"
$ws1.Range("D109").Value = 29
$ws1.Range("A110").Value = "08:48:01"
$ws1.Range("D110").Value = 33
$ws1.Range("A112").Value = "08:48:01"
$ws1.Range("D112").Value = 35
$ws1.Range("A115").Value = "08:48:01"
$ws1.Range("D115").Value = 36
$ws1.Range("A116").Value = "08:48:01"
$ws1.Range("D116").Value = 44
$ws1.Range("A117").Value = "08:48:01"
$ws1.Range("D117").Value = 45
$ws1.Range("A119").Value = "08:48:01"
$ws1.Range("C119").Value = "16_SANTA ANA"
$ws1.Range("D119").Value = 47
$ws1.Range("A120").Value = "08:48:01"
$ws1.Range("C120").Value = "23_HERNANDEZ"
$ws1.Range("D120").Value = 47
$ws1.Range("A122").Value = "08:48:01"
$ws1.Range("D122").Value = 54
$ws1.Range("A124").Value = "08:48:01"
$ws1.Range("D124").Value = 56
$ws1.Range("A125").Value = "08:48:01"
$ws1.Range("B125").Value = "09:52"
$ws1.Range("C125").Value = "15_ABASTO"
$ws1.Range("D125").Value = 64
$ws1.Range("A126").Value = "08:13:38"
$ws1.Range("B126").Value = "09:58"
$ws1.Range("C126").Value = "16_SANTA ANA"
$ws1.Range("D126").Value = 105
$ws1.Range("A127").Value = "08:48:01"
$ws1.Range("B127").Value = "10:11"
$ws1.Range("C127").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D127").Value = 83
$ws1.Range("B128").Value = "10:12"
$ws1.Range("C128").Value = "15_ABASTO"
$ws1.Range("D128").Value = 98
$ws1.Range("A129").Value = "08:48:01"
$ws1.Range("B129").Value = "10:21"
$ws1.Range("C129").Value = "26_HERNANDEZ"
$ws1.Range("D129").Value = 93
$ws1.Range("A130").Value = "08:48:01"
$ws1.Range("B130").Value = "10:22"
$ws1.Range("C130").Value = "17_ROMERO"
$ws1.Range("D130").Value = 94
$ws1.Range("A131").Value = "08:48:01"
$ws1.Range("B131").Value = "10:27"
$ws1.Range("C131").Value = "215A_EL PATO"
$ws1.Range("D131").Value = 99
$ws1.Range("E131").Value = "LP1912"
$ws1.Range("A132").Value = "08:48:01"
$ws1.Range("B132").Value = "10:42"
$ws1.Range("C132").Value = "17_ROMERO"
$ws1.Range("D132").Value = 114
$ws1.Range("E132").Value = "LP1912"
$ws1.Range("A133").Value = "08:48:01"
$ws1.Range("B133").Value = "10:44"
$ws1.Range("C133").Value = "14_ABASTO"
$ws1.Range("D133").Value = 116
$ws1.Range("E133").Value = "LP1912"

# ----- Sheet2: LP1912-215 -----
$ws2.Range("A2").Value = "Última actualización: 08:48:01"
$ws2.Range("A23").Value = "08:48:01"
$ws2.Range("D23").Value = 14
$ws2.Range("A25").Value = "08:48:01"
$ws2.Range("D25").Value = 54
$ws2.Range("A26").Value = "08:48:01"
$ws2.Range("D26").Value = 99

# ----- Sheet3: 6203-6173 -----
$ws3.Range("A2").Value = "Última actualización: 08:48:01"
$ws3.Range("A3").Value = "Total filas: 25"
$ws3.Range("A29").Value = "08:48:01"
$ws3.Range("B29").Value = "09:10"
$ws3.Range("C29").Value = "215D_LA PLATA"
$ws3.Range("D29").Value = 22
$ws3.Range("E29").Value = "L6203"
$ws3.Range("A30").Value = "08:48:01"
$ws3.Range("B30").Value = "10:03"
$ws3.Range("C30").Value = "215B_LP-P MOR-40 Y 115"
$ws3.Range("D30").Value = 75
$ws3.Range("E30").Value = "L6173"
